# Apply weekly update to "Hortaliza, Macroferia Regional de Talca - Apio" sheet.
# Two new price records are inserted into the weekly series:
#   - one for date 44720 (2022-06-08) at row 59
#   - one for date 44721 (2022-06-09) at (the then-current) row 75
# Inserting rows pushes all following records down accordingly, and the
# two previously-last records end up at the new final rows 183 and 184.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert first new row at position 59 ---------------------------------
$ws.Rows.Item(59).Insert()

$ws.Range("A59").Value2 = 5
$ws.Range("B59").Value2 = "Macroferia Regional de Talca"
$ws.Range("C59").Value2 = "Maule"
$ws.Range("D59").Value2 = 44720
$ws.Range("E59").Value2 = 7
$ws.Range("F59").Value2 = 100112017
$ws.Range("G59").Value2 = "Apio"
$ws.Range("H59").Value2 = "Americana (o)"
$ws.Range("I59").Value2 = "Primera"
$ws.Range("J59").Value2 = 700
$ws.Range("K59").Value2 = 6000
$ws.Range("L59").Value2 = 6000
$ws.Range("M59").Value2 = 6000
$ws.Range("N59").Value2 = "`$/docena de matas"
$ws.Range("O59").Value2 = "Provincia del Elquí"
$ws.Range("P59").Value2 = 1000
$ws.Range("Q59").Value2 = 6
$ws.Range("R59").Value2 = "Hortaliza"

# --- Insert second new row at position 75 (post first insert) ------------
$ws.Rows.Item(75).Insert()

$ws.Range("A75").Value2 = 5
$ws.Range("B75").Value2 = "Macroferia Regional de Talca"
$ws.Range("C75").Value2 = "Maule"
$ws.Range("D75").Value2 = 44721
$ws.Range("E75").Value2 = 7
$ws.Range("F75").Value2 = 100112017
$ws.Range("G75").Value2 = "Apio"
$ws.Range("H75").Value2 = "Americana (o)"
$ws.Range("I75").Value2 = "Primera"
$ws.Range("J75").Value2 = 500
$ws.Range("K75").Value2 = 6000
$ws.Range("L75").Value2 = 6000
$ws.Range("M75").Value2 = 6000
$ws.Range("N75").Value2 = "`$/docena de matas"
$ws.Range("O75").Value2 = "Provincia del Elquí"
$ws.Range("P75").Value2 = 1000
$ws.Range("Q75").Value2 = 6
$ws.Range("R75").Value2 = "Hortaliza"
